$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 additions ---
$ws.Range("F2").Value = "Patrz Panel administratora"
$ws.Range("G2").Value = "H"
$ws.Range("I2").Value = "datapocz i datakoniec zły format ma być char 10"

# --- New row 8 ---
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Panel administratora "
$ws.Range("C8").Value = "Kontrola importu"
$ws.Range("D8").Value = "Wordy\Panel admina importy.docx"
$ws.Range("G8").Value = "L"
$ws.Range("H8").Value = 43100

# reuse existing date style (style index 1) from H6
$ws.Range("H6").Copy()
$ws.Range("H8").PasteSpecial(-4122)

# add hyperlink on D8, then reuse hyperlink style from D2
$ws.Hyperlinks.Add($ws.Range("D8"), "Wordy\Panel%20admina%20importy.docx")
$ws.Range("D8").Value = "Wordy\Panel admina importy.docx"
$ws.Range("D2").Copy()
$ws.Range("D8").PasteSpecial(-4122)

$excel.CutCopyMode = 0

$ws.Range("D8").Select()
